$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "null" values in column G for the rows that held it
$gRows = @(2,3,4,6,8,10,11,12,14,15,18,19)
foreach ($r in $gRows) {
    $ws.Cells.Item($r, 7).Value = ""
}

# Update the sheet view / selection
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("B25").Select()
